$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Replace the two numeric cells on Sheet2 with text values. Order matters
# because it determines the order new entries are appended to the shared
# string table (G8 must become index 27, E5 must become index 28).
$ws2.Range("G8").Value = "0.345x"
$ws2.Range("E5").Value = "2000/6/14x"

# Update cell selections on each sheet, and make Sheet2 the active
# (tab-selected) sheet. Select Sheet1 and Sheet3 first, then select
# Sheet2 last so it ends up as the active tab.
$ws1.Range("B8").Select()
$ws3.Range("C28").Select()
$ws2.Range("K9").Select()
